# Make revenue-model template sheets business-agnostic: rename sheets to
# generic numbered "収益モデルN" slots and generalize their row labels /
# sample figures, plus clean up a couple of business-specific references
# on the assumptions sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the three business-specific sheets to generic numbered slots.
#    Renaming via COM automatically rewrites any formulas elsewhere in the
#    workbook (e.g. on "PL設計") that reference these sheet names.
# ---------------------------------------------------------------------
$wsModel1 = $wb.Worksheets.Item(2)   # formerly "ミールモデル"
$wsModel2 = $wb.Worksheets.Item(3)   # formerly "アカデミーモデル"
$wsModel3 = $wb.Worksheets.Item(4)   # formerly "コンサルモデル"

$wsModel1.Name = "収益モデル1"
$wsModel2.Name = "収益モデル2"
$wsModel3.Name = "収益モデル3"

# ---------------------------------------------------------------------
# 2) 収益モデル1 (formerly ミールモデル ユニットエコノミクス)
# ---------------------------------------------------------------------
$wsModel1.Range("A1").Value = "収益モデル1（セグメント1）"

$wsModel1.Range("A3").Value = "顧客数/取引数"
$wsModel1.Range("D3").Value = 360
$wsModel1.Range("E3").Value = 560
$wsModel1.Range("F3").Value = 800

$wsModel1.Range("A4").Value = "単価（円）"
$wsModel1.Range("D4").Value = 5500
$wsModel1.Range("E4").Value = 5500
$wsModel1.Range("F4").Value = 6000

$wsModel1.Range("A5").Value = "頻度/回数（月間）"
$wsModel1.Range("B5").Value = 3
$wsModel1.Range("C5").Value = 3

$wsModel1.Range("A6").Value = "成長率/解約率"

# Remove the LTV row (row 9) entirely and replace it with a thin spacer
# row, matching the pattern used elsewhere in the workbook (e.g. row 3 on
# the PL設計 sheet).
$wsModel1.Rows.Item(9).Clear()
$wsModel1.Rows.Item(9).RowHeight = 6

# ---------------------------------------------------------------------
# 3) 収益モデル2 (formerly アカデミーモデル 教育事業)
# ---------------------------------------------------------------------
$wsModel2.Range("A1").Value = "収益モデル2（セグメント2）"

$wsModel2.Range("A3").Value = "顧客数/取引数"
$wsModel2.Range("B3").Value = 75
$wsModel2.Range("C3").Value = 150
$wsModel2.Range("D3").Value = 270
$wsModel2.Range("E3").Value = 420
$wsModel2.Range("F3").Value = 600

$wsModel2.Range("A4").Value = "単価（円）"
$wsModel2.Range("B4").Value = 10000
$wsModel2.Range("C4").Value = 10000
$wsModel2.Range("D4").Value = 10500
$wsModel2.Range("E4").Value = 10500
$wsModel2.Range("F4").Value = 11000

$wsModel2.Range("A5").Value = "頻度/回数（月間）"
$wsModel2.Range("B5").Value = 3
$wsModel2.Range("C5").Value = 3
$wsModel2.Range("D5").Value = 4
$wsModel2.Range("E5").Value = 4

# Row 6 switches from an absolute "コンテンツ数" count to a
# "成長率/解約率" percentage, so set both value and number format to
# match the percentage style used by the equivalent row on the other
# model sheets.
$wsModel2.Range("A6").Value = "成長率/解約率"
$wsModel2.Range("B6").Value = 0.03
$wsModel2.Range("C6").Value = 0.03
$wsModel2.Range("D6").Value = 0.025
$wsModel2.Range("E6").Value = 0.025
$wsModel2.Range("F6").Value = 0.02
$wsModel2.Range("B6:F6").NumberFormat = "0.0%"

# Insert the thin spacer row 9 (didn't previously exist on this sheet).
$wsModel2.Rows.Item(9).RowHeight = 6

# ---------------------------------------------------------------------
# 4) 収益モデル3 (formerly コンサルモデル フライホイール)
# ---------------------------------------------------------------------
$wsModel3.Columns.Item(1).ColumnWidth = 23.1666666667

$wsModel3.Range("A1").Value = "収益モデル3（セグメント3）"

$wsModel3.Range("A3").Value = "顧客数/取引数"
$wsModel3.Range("B3").Value = 50
$wsModel3.Range("C3").Value = 100
$wsModel3.Range("D3").Value = 180
$wsModel3.Range("E3").Value = 280
$wsModel3.Range("F3").Value = 400

$wsModel3.Range("A4").Value = "単価（円）"
$wsModel3.Range("B4").Value = 15000
$wsModel3.Range("C4").Value = 15000
$wsModel3.Range("D4").Value = 15500
$wsModel3.Range("E4").Value = 15500
$wsModel3.Range("F4").Value = 16000

$wsModel3.Range("A5").Value = "頻度/回数（月間）"
$wsModel3.Range("B5").Value = 3
$wsModel3.Range("C5").Value = 3
$wsModel3.Range("D5").Value = 4
$wsModel3.Range("E5").Value = 4
$wsModel3.Range("F5").Value = 4

$wsModel3.Range("A6").Value = "成長率/解約率"
$wsModel3.Range("B6").Value = 0.03
$wsModel3.Range("C6").Value = 0.03
$wsModel3.Range("D6").Value = 0.025
$wsModel3.Range("E6").Value = 0.025
$wsModel3.Range("F6").Value = 0.02

# 月次売上 formula simplifies from a 4-factor product (incl. 稼働率) to a
# simple 顧客数/取引数 * 単価 product, matching the other two model sheets.
$wsModel3.Range("B7").Formula = "=B3*B4"
$wsModel3.Range("C7").Formula = "=C3*C4"
$wsModel3.Range("D7").Formula = "=D3*D4"
$wsModel3.Range("E7").Formula = "=E3*E4"
$wsModel3.Range("F7").Formula = "=F3*F4"

# Insert the thin spacer row 9 (didn't previously exist on this sheet).
$wsModel3.Rows.Item(9).RowHeight = 6

# ---------------------------------------------------------------------
# 5) （全Ver）前提条件: remove the two remaining business-specific notes.
# ---------------------------------------------------------------------
$wsAssumptions = $wb.Worksheets.Item(8)
$wsAssumptions.Range("C3").Value = "対象市場全体"
$wsAssumptions.Range("C8").Value = "変動費ベース"
